$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.784.55"
$ws.Range("E2").Value = "  -2.77%  "

$ws.Range("D3").Value = "2.745.78"
$ws.Range("E3").Value = "  -2.05%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "349.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.45%  "

$ws.Range("E7").Value = "  -3.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -3.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.68%  "

$ws.Range("E11").Value = "  +3.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0828"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.37%  "

$ws.Range("D15").Value = "3.176.43"
$ws.Range("E15").Value = "  -1.91%  "

$ws.Range("D16").Value = "2.742.95"
$ws.Range("E16").Value = "  -2.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.922"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").Value = "50.762.88"
$ws.Range("E18").Value = "  -2.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("E20").Value = "  -3.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.24%  "

$ws.Range("D22").Value = "0.0₃0952"
$ws.Range("E22").Value = "  -3.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.162"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.02%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.92%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "51.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.51%  "

$ws.Range("E32").Value = "  +3.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0441"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("E36").Value = "  -2.10%  "

$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.65%  "

$ws.Range("E40").Value = "  -4.01%  "

$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.14%  "

$ws.Range("E44").Value = "  -2.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.66%  "

$ws.Range("D46").Value = "2.076.25"
$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("E47").Value = "  -2.22%  "

$ws.Range("E48").Value = "  -1.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.905"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.44%  "

$ws.Range("E51").Value = "  +3.59%  "
